$d = $word.ActiveDocument

# --- 1. First paragraph: apply the "Titolo1" (Heading 1) style instead of
#        direct bold/size-30 character formatting. ---
$p1 = $d.Paragraphs(1)
$p1.Range.Style = "Titolo1"

# --- 2. Second paragraph: insert the missing word "mano " before "destra
#        agisce di conserva" and relocate the "_GoBack" bookmark to sit
#        right after the newly inserted word. ---
$p2 = $d.Paragraphs(2)
$r = $p2.Range.Duplicate()
$r.Find.Execute("La destra agisce", $true, $false, $false, $false, $false, $true, 1, $false, "La mano destra agisce", 2) | Out-Null

# $r now spans the replaced text "La mano destra agisce" -- use a temporary
# bookmark to force a run split right after "La " (before "mano "), so that
# "mano " ends up as its own run, matching a freshly-typed insertion.
$splitPos = $d.Range($r.Start + 3, $r.Start + 3)
$d.Bookmarks.Add("ZZZTempSplit", $splitPos)

# Relocate "_GoBack" to sit right after "mano " (i.e. right before "destra").
# Word keeps this bookmark unique/hidden and automatically drops it from its
# old location (end of paragraph 4) when it is re-added here.
$bmPos = $d.Range($r.Start + 8, $r.Start + 8)
$d.Bookmarks.Add("_GoBack", $bmPos)

# Remove the temporary bookmark, leaving the run split in place.
$d.Bookmarks("ZZZTempSplit").Delete()
